$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.366.74"
$ws.Range("E2").Value = "  -0.89%  "

$ws.Range("D3").Value = "2.642.48"
$ws.Range("E3").Value = "  -1.90%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.17%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "518.18"
$ws.Range("E5").Value = "  -0.48%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.67"
$ws.Range("E6").Value = "  -2.08%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  +0.36%  "

$ws.Range("E8").Value = "  -0.30%  "

$ws.Range("D9").Value = "2.653.50"
$ws.Range("E9").Value = "  -1.39%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.30"
$ws.Range("E10").Value = "  -3.44%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.105"
$ws.Range("E11").Value = "  -2.45%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.337"
$ws.Range("E12").Value = "  -1.41%  "

$ws.Range("E13").Value = "  +0.68%  "

$ws.Range("D14").Value = "3.108.00"
$ws.Range("E14").Value = "  -0.53%  "

$ws.Range("D15").Value = "59.362.40"
$ws.Range("E15").Value = "  -0.68%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "20.94"
$ws.Range("E16").Value = "  -2.74%  "

$ws.Range("E17").Value = "  -1.65%  "

$ws.Range("D18").Value = "2.655.54"
$ws.Range("E18").Value = "  -2.37%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "349.18"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.50"
$ws.Range("E20").Value = "  -3.34%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.34"
$ws.Range("E21").Value = "  -2.80%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.21"
$ws.Range("E22").Value = "  -0.92%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  +0.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.99"
$ws.Range("E24").Value = "  +1.18%  "

$ws.Range("E25").Value = "  -2.64%  "

$ws.Range("E26").Value = "  +1.82%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  +0.40%  "

$ws.Range("D28").Value = "0.0₃0805"
$ws.Range("E28").Value = "  -3.38%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.14"
$ws.Range("E29").Value = "  -1.58%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.25%  "

$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.92"
$ws.Range("E32").Value = "  -1.27%  "

$ws.Range("B33").Value = "PancakeSwap"
$ws.Range("C33").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.58"
$ws.Range("E33").Value = "  -1.03%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "149.69"
$ws.Range("E34").Value = "  +0.66%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.06"
$ws.Range("E35").Value = "  -0.53%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.946"
$ws.Range("E36").Value = "  -12.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.19"
$ws.Range("E37").Value = "  +1.51%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.862"
$ws.Range("E38").Value = "  -2.13%  "

$ws.Range("E39").Value = "  +0.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.48"
$ws.Range("E40").Value = "  +2.34%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.68"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0994"
$ws.Range("E42").Value = "  -0.37%  "

$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  +0.02%  "

$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "276.68"
$ws.Range("E44").Value = "  -3.66%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.602"
$ws.Range("E45").Value = "  -3.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.64"
$ws.Range("E46").Value = "  -1.31%  "

$ws.Range("D47").Value = "2.089.38"
$ws.Range("E47").Value = "  +4.52%  "

$ws.Range("E48").Value = "  -3.72%  "

$ws.Range("E49").Value = "  -0.92%  "

$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.31"
$ws.Range("E50").Value = "  +0.60%  "

$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.71"
$ws.Range("E51").Value = "  -1.88%  "
